# DataPruebas.xlsx - cierre de la segunda unidad del curso
#
# Adds the "CP002_iniciar_sesion" test case (row 3), which previously was
# only a placeholder "CP002" row, and tweaks the "CP001_creacion_cta" row
# (row 2) with a new e-mail/user value, highlights the title cell and moves
# the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataPruebas")

# --- Highlight the title cell (A1) with a yellow fill -----------------
$ws.Range("A1").Interior.Color = 65535

# --- Row 3: flesh out the "CP002_iniciar_sesion" test case ------------
# (was just a stub row with "CP002" / "X" placeholder values)
$ws.Range("A3").Value = "CP002_iniciar_sesion"
$ws.Range("B3").Value = "e4587y87678@algo.com"
$ws.Range("C3").Value = "adasdas12"
$ws.Range("D3").Value = "USer Test 001"

# --- Row 2: update the "CP001_creacion_cta" test case ------------------
$ws.Range("D2").Value = "Pobre Dgo"
$ws.Range("I2").Value = "Pobre Dgo"
$ws.Range("B2").Value = "e4r5dfs432fs8@algo.com"

# --- Hyperlinks: rebuild so B3's new mail link becomes rId1 and B2's
#     (changed) address becomes rId2, in that order ---------------------
$ws.Range("B2").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:e4587y87678@algo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:e4r5dfs432fs8@algo.com") | Out-Null

# Reapply the shared "Hipervínculo" style so both hyperlink cells reuse
# the same style record instead of each getting their own.
$ws.Range("B2").Style = "Hipervínculo"
$ws.Range("B3").Style = "Hipervínculo"

# --- Move the active selection to C13 -----------------------------------
$ws.Activate()
$ws.Range("C13").Select() | Out-Null
